{"js": "// Rename the 2nd and 3rd header-row columns of the feature table:\n//   \"Customer Benefits\"   -> \"Customers Benefits\"\n//   \"Management Features\" -> \"Professionals Features\"\n// (column formatting - bold, Segoe UI, border, color - is preserved\n// because we only change the text value of the existing cell, not its\n// paragraph/run properties.)\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"No tables found in the document.\");\n}\n\nconst table = tables.items[0];\n\n// Header row is row index 0: [\"Feature\", \"Customer Benefits\", \"Management Features\"]\nconst customerBenefitsCell = table.getCell(0, 1);\nconst managementFeaturesCell = table.getCell(0, 2);\n\ncustomerBenefitsCell.load(\"value\");\nmanagementFeaturesCell.load(\"value\");\nawait context.sync();\n\nif (customerBenefitsCell.value === \"Customer Benefits\") {\n  customerBenefitsCell.value = \"Customers Benefits\";\n}\n\nif (managementFeaturesCell.value === \"Management Features\") {\n  managementFeaturesCell.value = \"Professionals Features\";\n}\n\nawait context.sync();\n", "ps1": "# Rename the 2nd and 3rd header-row columns of the feature table:\n#   \"Customer Benefits\"   -> \"Customers Benefits\"\n#   \"Management Features\" -> \"Professionals Features\"\n# Only the cell text is touched, so the existing run formatting\n# (bold, Segoe UI, border, color) carries over untouched.\n\n$d = $word.ActiveDocument\n$table = $d.Tables.Item(1)\n\n$headerRow = 1\nfor ($col = 1; $col -le $table.Columns.Count; $col++) {\n    $cell = $table.Cell($headerRow, $col)\n    $cellText = $cell.Range.Text\n    $cellText = $cellText.TrimEnd([char]7, [char]13)\n\n    if ($cellText -eq \"Customer Benefits\") {\n        $cell.Range.Text = \"Customers Benefits\"\n    }\n    elseif ($cellText -eq \"Management Features\") {\n        $cell.Range.Text = \"Professionals Features\"\n    }\n}\n"}
